# Fix SMS sending (kavenegar) data: correct a phone number and rotate
# the generated passwords for the two users in the sheet, then tidy up
# the column widths so the new phone numbers / passwords are fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 ("علی عمادی") had a wrong phone number recorded - fix it.
$ws.Range("D2").Value = 9381088235

# Rotate the (placeholder) passwords for both users.
$ws.Range("F2").Value = "4kTWgnoA"
$ws.Range("F3").Value = "7Ab1160p"

# Re-size the phone & password columns so the new values fit nicely.
$ws.Columns.Item(4).ColumnWidth = 10.8333333333
$ws.Columns.Item(6).ColumnWidth = 14.1666666666

# Leave the selection where the user last clicked.
$ws.Range("K11").Select()

$wb.Save()
